# Updates the Price (D) and Volume(1h) (E) columns of the cryptos table
# to the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = <new price text>; E = <new volume text> }
# D is omitted for rows whose price did not change; same for E.
$updates = [ordered]@{
    2 = @{ D='57.692.18'; E='  -1.72%  ' }
    3 = @{ D='2.275.44'; E='  -1.26%  ' }
    4 = @{ E='  +0.14%  ' }
    5 = @{ D='530.81'; E='  -3.48%  ' }
    6 = @{ D='130.90'; E='  -0.08%  ' }
    7 = @{ E='  +0.16%  ' }
    8 = @{ D='0.580'; E='  +1.94%  ' }
    9 = @{ D='2.274.15'; E='  -1.48%  ' }
    10 = @{ D='0.0989'; E='  -3.35%  ' }
    11 = @{ D='5.45'; E='  -1.80%  ' }
    12 = @{ E='  +0.26%  ' }
    13 = @{ D='0.329'; E='  -2.27%  ' }
    14 = @{ D='23.27'; E='  -2.17%  ' }
    15 = @{ D='2.689.57'; E='  -1.15%  ' }
    16 = @{ D='57.719.07'; E='  -1.65%  ' }
    17 = @{ D='0.0000131'; E='  -1.76%  ' }
    18 = @{ D='2.247.85'; E='  -2.34%  ' }
    19 = @{ D='10.49'; E='  -2.67%  ' }
    20 = @{ D='4.15'; E='  -5.07%  ' }
    21 = @{ D='310.37'; E='  -1.62%  ' }
    22 = @{ D='6.36'; E='  -2.42%  ' }
    23 = @{ E='  -0.11%  ' }
    24 = @{ D='62.28'; E='  -1.30%  ' }
    25 = @{ D='0.166'; E='  -2.70%  ' }
    26 = @{ D='1.00'; E='  +0.21%  ' }
    27 = @{ D='7.89'; E='  -4.05%  ' }
    28 = @{ D='1.24'; E='  -5.76%  ' }
    29 = @{ D='169.24'; E='  -0.90%  ' }
    30 = @{ D='1.69'; E='  -4.16%  ' }
    31 = @{ D='0.0₃0715'; E='  -2.36%  ' }
    32 = @{ D='5.70'; E='  -2.77%  ' }
    33 = @{ D='1.03'; E='  -4.99%  ' }
    34 = @{ D='0.376'; E='  -3.16%  ' }
    35 = @{ E='  -0.03%  ' }
    36 = @{ D='17.67'; E='  -0.56%  ' }
    37 = @{ E='  -0.11%  ' }
    38 = @{ D='1.22'; E='  -4.25%  ' }
    39 = @{ D='3.84'; E='  -3.92%  ' }
    40 = @{ D='38.39'; E='  -0.04%  ' }
    41 = @{ D='1.47'; E='  -3.17%  ' }
    42 = @{ D='138.03'; E='  -2.77%  ' }
    43 = @{ D='283.45'; E='  -4.45%  ' }
    44 = @{ D='3.39'; E='  -2.04%  ' }
    45 = @{ D='0.0941'; E='  -0.55%  ' }
    46 = @{ D='0.0491'; E='  -1.67%  ' }
    47 = @{ D='0.549'; E='  -1.24%  ' }
    48 = @{ D='17.81'; E='  -4.38%  ' }
    49 = @{ D='0.0208'; E='  -2.86%  ' }
    50 = @{ D='10.93' }
    51 = @{ E='  -0.57%  ' }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    if ($rowData.ContainsKey("D")) {
        # Prefix with a leading apostrophe so Excel stores the price as literal
        # text (matching the source inline-string cell) instead of re-parsing
        # numeric-looking text (e.g. "530.81", multi-dot "57.692.18") into a
        # floating point number. Reset the style afterwards so the quote-prefix
        # formatting does not leave a stray style on the cell.
        $cell = $ws.Range("D" + $row)
        $cell.Value = "'" + $rowData["D"]
        $cell.Style = "Normal"
    }
    if ($rowData.ContainsKey("E")) {
        $ws.Range("E" + $row).Value = $rowData["E"]
    }
}
